$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "SAP Regression Automation"
$ws.Range("B5").Value = "Yes"
$ws.Range("C5").Value = "AU_OC_04"
$ws.Range("D5").Value = 5400338

$ws.Range("A6").Value = "SAP Regression Automation"
$ws.Range("B6").Value = "Yes"
$ws.Range("C6").Value = "AU_OC_05"
$ws.Range("D6").Value = 5400339

$ws.Range("A7").Value = "SAP Regression Automation"
$ws.Range("B7").Value = "Yes"
$ws.Range("C7").Value = "AU_OC_06"
$ws.Range("D7").Value = 5400340

$ws.Range("E7").Value = "TC_06 Trigger Invoice Order with Existing Customer"
$ws.Range("E6").Value = "TC_05 Trigger CreditCard Order with Existing Customer"
$ws.Range("E5").Value = "TC_04 Trigger Alipay Order with Existing Customer"
